$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new reading was recorded for 2026/01/19 at 16:00 (value 122) that sorts
# between the existing 13:00 row (681) and the 2026/12/29 block that used to
# start at row 682. Insert a fresh row so everything below shifts down by one,
# then populate it.
$ws.Rows("682:682").Insert()

# Column A holds dates stored as plain text (e.g. "2026/01/19"), not real
# date values. Assigning a date-looking string directly makes Excel's
# autodetection coerce it into a date serial + date number format, so we
# temporarily force Text formatting for the assignment and then clear the
# formatting again to fall back to the sheet's default (unstyled) cell,
# matching every other row in this column.
$ws.Range("A682").NumberFormat = "@"
$ws.Range("A682").Value = "2026/01/19"
$ws.Range("A682").ClearFormats()

$ws.Range("B682").Value = "月"
$ws.Range("C682").Value = 16
$ws.Range("D682").Value = 122
